$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-14 from 45221 (2023-10-22) to 45224 (2023-10-25)
$ws.Range("C2:C14").Value = 45224
